$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "52.145.56"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.907.66"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +3.69%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "350.88"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "112.18"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.557"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.91"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("E11").Value = "  +2.55%  "
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.93"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.79"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.364.37"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.75%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.913.50"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.08%  "
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.993"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +5.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "52.180.48"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.62"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.30"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.15"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0978"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.77"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.53"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.79"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.67"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.01%  "
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.57"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.59"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.25"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.44"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.06"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0957"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +11.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "53.06"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0451"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.07"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.40%  "
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.64"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.34%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.82"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +13.42%  "
$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.08"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.84%  "
$ws.Range("E42").Value = "  +1.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.36"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +6.17%  "
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.64"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +7.00%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "120.99"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.67%  "
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.197.94"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +4.29%  "
$ws.Range("E48").Value = "  +3.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.260"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +22.00%  "
$ws.Range("B50").Value = "SEI"
$ws.Range("C50").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.960"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.38%  "
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0334"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +11.54%  "
